$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cells that no longer exist in the updated data (sales_invested_capital,
# roic, roic_cost_capital columns were dropped for these rows)
$ws.Range("Z2:Z3").ClearContents()
$ws.Range("AA2:AA3").ClearContents()
$ws.Range("AC2:AC3").ClearContents()

# Row 2 (first Netherlands Apparel company) updated metrics
$ws.Range("G2").Value = -0.8673965936739658
$ws.Range("H2").Value = -0.8673965936739658
$ws.Range("I2").Value = -0.8917274939172749
$ws.Range("J2").Value = -0.8917274939172749
$ws.Range("K2").Value = -96.59999999999999
$ws.Range("L2").Value = -1.175182481751825
$ws.Range("O2").Value = -0.0
$ws.Range("R2").Value = -0.0
$ws.Range("U2").Value = 6.6
$ws.Range("V2").Value = 0.01343648208469055
$ws.Range("X2").Value = 0.06702908493975343
$ws.Range("AB2").Value = 0.05073555458025421
$ws.Range("AD2").Value = 283.3
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 283.3
$ws.Range("AG2").Value = 276.7
$ws.Range("AH2").Value = 0.3657843770174306
$ws.Range("AI2").Value = 2.384680134680135
$ws.Range("AJ2").Value = 0.3603333767417632
$ws.Range("AK2").Value = 2.466131907308378
$ws.Range("AL2").Value = 20.2
$ws.Range("AM2").Value = 20.175
$ws.Range("AN2").Value = -4.234678624813154
$ws.Range("AO2").Value = -3.628712871287129
$ws.Range("AP2").Value = -4.136023916292974
$ws.Range("AQ2").Value = -3.633209417596035

# Row 3 (La Perla Fashion Holding N.V.) updated metrics
$ws.Range("G3").Value = -0.8673965936739658
$ws.Range("H3").Value = -0.8673965936739658
$ws.Range("I3").Value = -0.8917274939172749
$ws.Range("J3").Value = -0.8917274939172749
$ws.Range("K3").Value = -96.59999999999999
$ws.Range("L3").Value = -1.175182481751825
$ws.Range("O3").Value = 0.0
$ws.Range("R3").Value = 0.0
$ws.Range("U3").Value = 6.6
$ws.Range("V3").Value = 0.01343648208469055
$ws.Range("X3").Value = 0.06702908493975343
$ws.Range("AB3").Value = 0.05073555458025421
$ws.Range("AD3").Value = 283.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 283.3
$ws.Range("AG3").Value = 276.7
$ws.Range("AH3").Value = 0.3657843770174306
$ws.Range("AI3").Value = 2.384680134680135
$ws.Range("AJ3").Value = 0.3603333767417632
$ws.Range("AK3").Value = 2.466131907308378
$ws.Range("AL3").Value = 20.2
$ws.Range("AM3").Value = 20.175
$ws.Range("AN3").Value = -4.234678624813154
$ws.Range("AO3").Value = -3.628712871287129
$ws.Range("AP3").Value = -4.136023916292974
$ws.Range("AQ3").Value = -3.633209417596035
